$wb = $excel.ActiveWorkbook

# --- Sheet: x_df ---
$ws = $wb.Worksheets.Item("x_df")
$ws.Range("B4").Value = 97.64
$ws.Range("B5").Value = 93.75
$ws.Range("B6").Value = 93.37
$ws.Range("B7").Value = 93.33
$ws.Range("B8").Value = 93.33
$ws.Range("B11").Value = 89.56

# --- Sheet: q_df ---
$ws = $wb.Worksheets.Item("q_df")
$ws.Range("A2").Value = 81.93936947794616
$ws.Range("B2").Value = 73.45496719657494
$ws.Range("C2").Value = 71.46114056540813
$ws.Range("D2").Value = 76.56831891577663
$ws.Range("E2").Value = 78.50307698478537
$ws.Range("F2").Value = 80.57878795051627
$ws.Range("G2").Value = 82.01886235588172
$ws.Range("H2").Value = 83.49893334848616
$ws.Range("I2").Value = 85.23838890122474
$ws.Range("J2").Value = 86.74303483638201
$ws.Range("K2").Value = 91.88970435195002
$ws.Range("A3").Value = 81.93936947794616
$ws.Range("B3").Value = 73.45496719657494
$ws.Range("C3").Value = 71.46114056540813
$ws.Range("D3").Value = 76.56831891577663
$ws.Range("E3").Value = 78.50307698478537
$ws.Range("F3").Value = 80.57878795051627
$ws.Range("G3").Value = 82.01886235588172
$ws.Range("H3").Value = 83.49893334848616
$ws.Range("I3").Value = 85.23838890122474
$ws.Range("J3").Value = 86.74303483638201
$ws.Range("K3").Value = 91.88970435195002
$ws.Range("A4").Value = 81.83194065223961
$ws.Range("B4").Value = 74.74535188773635
$ws.Range("C4").Value = 72.88202405734792
$ws.Range("D4").Value = 77.27524520582892
$ws.Range("E4").Value = 78.62869189895402
$ws.Range("F4").Value = 80.16143279691184
$ws.Range("G4").Value = 81.52676427473973
$ws.Range("H4").Value = 83.0283825776067
$ws.Range("I4").Value = 85.07037489744788
$ws.Range("J4").Value = 87.84418794740816
$ws.Range("K4").Value = 94.09990444760692
$ws.Range("A5").Value = 81.65528076353313
$ws.Range("B5").Value = 75.78554483863378
$ws.Range("C5").Value = 74.78243767904854
$ws.Range("D5").Value = 77.40078808070729
$ws.Range("E5").Value = 78.10066444499289
$ws.Range("F5").Value = 79.29567576097303
$ws.Range("G5").Value = 80.56187224731259
$ws.Range("H5").Value = 82.50939666541798
$ws.Range("I5").Value = 86.12790627262274
$ws.Range("J5").Value = 90.81349408987724
$ws.Range("K5").Value = 99.14248434247484
$ws.Range("A6").Value = 81.63821313664103
$ws.Range("B6").Value = 75.82167835099672
$ws.Range("C6").Value = 74.73038858488631
$ws.Range("D6").Value = 77.33394100779675
$ws.Range("E6").Value = 78.00362865271136
$ws.Range("F6").Value = 79.20419001685895
$ws.Range("G6").Value = 80.49197263998364
$ws.Range("H6").Value = 82.47372967994528
$ws.Range("I6").Value = 86.17155893759588
$ws.Range("J6").Value = 91.05119266903736
$ws.Range("K6").Value = 99.75217943060038
$ws.Range("A7").Value = 81.63618631913042
$ws.Range("B7").Value = 75.82348488166237
$ws.Range("C7").Value = 74.7242076426776
$ws.Range("D7").Value = 77.33964343972544
$ws.Range("E7").Value = 77.9886129556204
$ws.Range("F7").Value = 79.19492661693559
$ws.Range("G7").Value = 80.48687982127764
$ws.Range("H7").Value = 82.48596428461953
$ws.Range("I7").Value = 86.18741869982271
$ws.Range("J7").Value = 91.12221978814586
$ws.Range("K7").Value = 99.82458203750988
$ws.Range("A8").Value = 81.63618631913042
$ws.Range("B8").Value = 75.82348488166237
$ws.Range("C8").Value = 74.7242076426776
$ws.Range("D8").Value = 77.33964343972544
$ws.Range("E8").Value = 77.9886129556204
$ws.Range("F8").Value = 79.19492661693559
$ws.Range("G8").Value = 80.48687982127764
$ws.Range("H8").Value = 82.48596428461953
$ws.Range("I8").Value = 86.18741869982271
$ws.Range("J8").Value = 91.12221978814586
$ws.Range("K8").Value = 99.82458203750988
$ws.Range("A9").Value = 81.63618631913042
$ws.Range("B9").Value = 75.82348488166237
$ws.Range("C9").Value = 74.7242076426776
$ws.Range("D9").Value = 77.33964343972544
$ws.Range("E9").Value = 77.9886129556204
$ws.Range("F9").Value = 79.19492661693559
$ws.Range("G9").Value = 80.48687982127764
$ws.Range("H9").Value = 82.48596428461953
$ws.Range("I9").Value = 86.18741869982271
$ws.Range("J9").Value = 91.12221978814586
$ws.Range("K9").Value = 99.82458203750988
$ws.Range("A10").Value = 81.63618631913042
$ws.Range("B10").Value = 75.82348488166237
$ws.Range("C10").Value = 74.7242076426776
$ws.Range("D10").Value = 77.33964343972544
$ws.Range("E10").Value = 77.9886129556204
$ws.Range("F10").Value = 79.19492661693559
$ws.Range("G10").Value = 80.48687982127764
$ws.Range("H10").Value = 82.48596428461953
$ws.Range("I10").Value = 86.18741869982271
$ws.Range("J10").Value = 91.12221978814586
$ws.Range("K10").Value = 99.82458203750988
$ws.Range("A11").Value = 81.46465583119641
$ws.Range("B11").Value = 75.5089449774794
$ws.Range("C11").Value = 74.98897253472948
$ws.Range("D11").Value = 76.49997115894458
$ws.Range("E11").Value = 77.18105468021393
$ws.Range("F11").Value = 78.28970801900307
$ws.Range("G11").Value = 79.63355254679018
$ws.Range("H11").Value = 82.02601548596317
$ws.Range("I11").Value = 87.69665586835072
$ws.Range("J11").Value = 94.50774743486019
$ws.Range("K11").Value = 106.53725810095804

# --- Sheet: q_df_Test ---
$ws = $wb.Worksheets.Item("q_df_Test")
$ws.Range("A2").Value = 81.99190297005694
$ws.Range("B2").Value = 74.0491760456403
$ws.Range("C2").Value = 72.37672751229083
$ws.Range("D2").Value = 76.75170523180327
$ws.Range("E2").Value = 78.42804310247209
$ws.Range("F2").Value = 80.45023469847035
$ws.Range("G2").Value = 82.06108680795897
$ws.Range("H2").Value = 83.46610486422708
$ws.Range("I2").Value = 85.12507615229177
$ws.Range("J2").Value = 86.45094560446282
$ws.Range("K2").Value = 93.97199436860471
$ws.Range("A3").Value = 81.99190297005694
$ws.Range("B3").Value = 74.0491760456403
$ws.Range("C3").Value = 72.37672751229083
$ws.Range("D3").Value = 76.75170523180327
$ws.Range("E3").Value = 78.42804310247209
$ws.Range("F3").Value = 80.45023469847035
$ws.Range("G3").Value = 82.06108680795897
$ws.Range("H3").Value = 83.46610486422708
$ws.Range("I3").Value = 85.12507615229177
$ws.Range("J3").Value = 86.45094560446282
$ws.Range("K3").Value = 93.97199436860471
$ws.Range("A4").Value = 81.86553135906756
$ws.Range("B4").Value = 75.08005066876817
$ws.Range("C4").Value = 73.13697853605835
$ws.Range("D4").Value = 77.23151682911471
$ws.Range("E4").Value = 78.42588711828934
$ws.Range("F4").Value = 80.07235649461062
$ws.Range("G4").Value = 81.54466296777733
$ws.Range("H4").Value = 83.09820299571483
$ws.Range("I4").Value = 85.0262743517023
$ws.Range("J4").Value = 87.38898997181445
$ws.Range("K4").Value = 96.74363932462612
$ws.Range("A5").Value = 81.65772126167327
$ws.Range("B5").Value = 75.9343997108673
$ws.Range("C5").Value = 74.61925652705791
$ws.Range("D5").Value = 77.26346633272786
$ws.Range("E5").Value = 78.02399578663815
$ws.Range("F5").Value = 79.37892144894387
$ws.Range("G5").Value = 80.81105224064686
$ws.Range("H5").Value = 82.44708139808967
$ws.Range("I5").Value = 85.24545873200816
$ws.Range("J5").Value = 89.73781589444846
$ws.Range("K5").Value = 100.83691246599874
$ws.Range("A6").Value = 81.63764412230279
$ws.Range("B6").Value = 75.9361794889809
$ws.Range("C6").Value = 74.59918828591381
$ws.Range("D6").Value = 77.20373929887188
$ws.Range("E6").Value = 77.9487532121571
$ws.Range("F6").Value = 79.34419794933474
$ws.Range("G6").Value = 80.72412042579957
$ws.Range("H6").Value = 82.43229111972138
$ws.Range("I6").Value = 85.2583522683476
$ws.Range("J6").Value = 89.79079207226889
$ws.Range("K6").Value = 101.18788121683912
$ws.Range("A7").Value = 81.63525991869342
$ws.Range("B7").Value = 75.93524385306277
$ws.Range("C7").Value = 74.59680513898803
$ws.Range("D7").Value = 77.19401877233565
$ws.Range("E7").Value = 77.93749362957232
$ws.Range("F7").Value = 79.3318882821361
$ws.Range("G7").Value = 80.70910472870861
$ws.Range("H7").Value = 82.42108047564709
$ws.Range("I7").Value = 85.25988340360101
$ws.Range("J7").Value = 89.79708310766182
$ws.Range("K7").Value = 101.22955951308842
$ws.Range("A8").Value = 81.63525991869342
$ws.Range("B8").Value = 75.93524385306277
$ws.Range("C8").Value = 74.59680513898803
$ws.Range("D8").Value = 77.19401877233565
$ws.Range("E8").Value = 77.93749362957232
$ws.Range("F8").Value = 79.3318882821361
$ws.Range("G8").Value = 80.70910472870861
$ws.Range("H8").Value = 82.42108047564709
$ws.Range("I8").Value = 85.25988340360101
$ws.Range("J8").Value = 89.79708310766182
$ws.Range("K8").Value = 101.22955951308842
$ws.Range("A9").Value = 81.63525991869342
$ws.Range("B9").Value = 75.93524385306277
$ws.Range("C9").Value = 74.59680513898803
$ws.Range("D9").Value = 77.19401877233565
$ws.Range("E9").Value = 77.93749362957232
$ws.Range("F9").Value = 79.3318882821361
$ws.Range("G9").Value = 80.70910472870861
$ws.Range("H9").Value = 82.42108047564709
$ws.Range("I9").Value = 85.25988340360101
$ws.Range("J9").Value = 89.79708310766182
$ws.Range("K9").Value = 101.22955951308842
$ws.Range("A10").Value = 81.63525991869342
$ws.Range("B10").Value = 75.93524385306277
$ws.Range("C10").Value = 74.59680513898803
$ws.Range("D10").Value = 77.19401877233565
$ws.Range("E10").Value = 77.93749362957232
$ws.Range("F10").Value = 79.3318882821361
$ws.Range("G10").Value = 80.70910472870861
$ws.Range("H10").Value = 82.42108047564709
$ws.Range("I10").Value = 85.25988340360101
$ws.Range("J10").Value = 89.79708310766182
$ws.Range("K10").Value = 101.22955951308842
$ws.Range("A11").Value = 81.4334836826465
$ws.Range("B11").Value = 75.65215050072545
$ws.Range("C11").Value = 75.20577313119506
$ws.Range("D11").Value = 76.48545647436838
$ws.Range("E11").Value = 77.19304281051501
$ws.Range("F11").Value = 78.36790803994197
$ws.Range("G11").Value = 79.78670226774847
$ws.Range("H11").Value = 82.02825435680636
$ws.Range("I11").Value = 86.39482739713175
$ws.Range("J11").Value = 91.68747430285129
$ws.Range("K11").Value = 108.00401266013469
